$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.692.85"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.776.49"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.27"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").Value = "3.772.76"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.443"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.59"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "4.411.12"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "3.752.12"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "67.718.27"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "456.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.689"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.82"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("E25").Value = "  -7.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "3.926.43"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -6.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.77"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.90"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0989"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  +4.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.77"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.79%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.292"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "383.30"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.17%  "
